$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 11 (Marking): B11 6 -> 9, C11 3 -> 2
$ws.Range("B11").Value = 9
$ws.Range("C11").Value = 2

# Row 12 (Total): B12 96 -> 144, C12 -24 -> -16, E12 "72/168" -> "128/252"
$ws.Range("B12").Value = 144
$ws.Range("C12").Value = -16
$ws.Range("E12").Value = "128/252"
